$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 97.90318578020005
$ws.Range("H2").Value = 97.62862158860304
$ws.Range("I2").Value = 96.43863793986142

$ws.Range("G3").Value = 98.17988368673714
$ws.Range("H3").Value = 97.63424170044945
$ws.Range("I3").Value = 96.52950974918801

$ws.Range("G4").Value = 98.02710530101348
$ws.Range("H4").Value = 97.58046813653365
$ws.Range("I4").Value = 96.19771780083322

$ws.Range("G5").Value = 97.95539475453438
$ws.Range("H5").Value = 97.59099172472406
$ws.Range("I5").Value = 96.2343149660089

$ws.Range("G6").Value = 98.13222834774639
$ws.Range("H6").Value = 97.60101959686401
$ws.Range("I6").Value = 96.05793635543462
